# Update countries & provincias Spain
# - Reorder Bulgaria/Jordania and Maldivas/Consejo Danes entries (data refresh
#   moved Bulgaria & Maldivas ahead of their former neighbours in the table,
#   which is sorted by total cases).
# - Refresh the "Datos actualizados" timestamp.
# - Refresh numeric statistics for several countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 23:47"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 7940007
$ws.Range("C4").Value = 45529
$ws.Range("D4").Value = 5083978
$ws.Range("E4").Value = 2636797
$ws.Range("G4").Value = 585
$ws.Range("H4").Value = 219232

# --- Row 6: Brasil ---
$ws.Range("B6").Value = 5082637
$ws.Range("C6").Value = 25447
$ws.Range("E6").Value = 498844
$ws.Range("G6").Value = 506
$ws.Range("H6").Value = 150198

# --- Rows 85/86: Bulgaria moves ahead of Jordania (both keep their own
#     numbers; Bulgaria also gets fresh totals) ---
$ws.Range("A85").Value = "Bulgaria"
$ws.Range("B85").Value = 24319
$ws.Range("C85").Value = 448
$ws.Range("D85").Value = 15818
$ws.Range("E85").Value = 7610
$ws.Range("G85").Value = 4
$ws.Range("H85").Value = 891

$ws.Range("A86").Value = "Jordania"
$ws.Range("B86").Value = 23998
$ws.Range("C86").Value = 1235
$ws.Range("D86").Value = 6045
$ws.Range("E86").Value = 17772
$ws.Range("G86").Value = 15
$ws.Range("H86").Value = 181

# --- Rows 105/106: Maldivas moves ahead of Consejo Danes para los
#     Refugiados (Consejo Danes keeps its former numbers; Maldivas gets
#     fresh totals) ---
$ws.Range("A105").Value = "Maldivas"
$ws.Range("B105").Value = 10859
$ws.Range("C105").Value = 51
$ws.Range("D105").Value = 9683
$ws.Range("E105").Value = 1142
$ws.Range("H105").Value = 34

$ws.Range("A106").Value = "Consejo Danes para los Refugiados"
$ws.Range("B106").Value = 10841
$ws.Range("C106").Value = 6
$ws.Range("D106").Value = 10242
$ws.Range("E106").Value = 323
$ws.Range("H106").Value = 276

# --- Row 108: Guayana Francesa ---
$ws.Range("D108").Value = 9834
$ws.Range("E108").Value = 241

# --- Row 118: Cabo Verde ---
$ws.Range("B118").Value = 6913
$ws.Range("C118").Value = 104
$ws.Range("D118").Value = 5970
$ws.Range("E118").Value = 869
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 74

# --- Row 130: Trinidad yTobago ---
$ws.Range("B130").Value = 5043
$ws.Range("C130").Value = 22
$ws.Range("D130").Value = 3221
$ws.Range("E130").Value = 1732
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 90

# --- Row 161: Togo ---
$ws.Range("B161").Value = 1935
$ws.Range("C161").Value = 14
$ws.Range("D161").Value = 1444
$ws.Range("E161").Value = 442
